# Applies the cryptos-list price/volume refresh described by the commit
# "Updated cryptos list on Wed May 29 21:14:34 UTC 2024 with GitHub Actions".
#
# Every touched cell in the source workbook is stored as literal text (no
# "s" style attribute, t="inlineStr"/shared-string) even though many of the
# Price values look like plain numbers (e.g. "6.47"). Writing such a string
# straight into Range.Value makes Excel auto-convert it to a real number,
# so we prefix with a literal leading apostrophe to force text entry (exactly
# like typing '6.47 into a cell) and then reset .Style back to "Normal" so the
# implicit quote-prefix formatting Excel applies doesn't leave a style-index
# footprint behind - matching the original (unstyled) cells exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $text
    $range.Style = 'Normal'
}

# Row 2
Set-TextValue 'D2' '67.411.29'
Set-TextValue 'E2' '  -1.32%  '

# Row 3
Set-TextValue 'D3' '3.751.87'
Set-TextValue 'E3' '  -2.11%  '

# Row 4
Set-TextValue 'E4' '  +0.06%  '

# Row 5
Set-TextValue 'D5' '594.59'
Set-TextValue 'E5' '  -1.13%  '

# Row 6
Set-TextValue 'D6' '168.41'
Set-TextValue 'E6' '  -0.62%  '

# Row 7
Set-TextValue 'D7' '3.749.71'
Set-TextValue 'E7' '  -2.14%  '

# Row 8
Set-TextValue 'E8' '  -0.19%  '

# Row 9
Set-TextValue 'D9' '0.521'
Set-TextValue 'E9' '  -1.09%  '

# Row 10
Set-TextValue 'D10' '0.163'
Set-TextValue 'E10' '  -1.47%  '

# Row 11
Set-TextValue 'D11' '6.47'
Set-TextValue 'E11' '  -0.22%  '

# Row 12
Set-TextValue 'D12' '0.451'
Set-TextValue 'E12' '  -1.44%  '

# Row 13
Set-TextValue 'D13' '0.0000272'
Set-TextValue 'E13' '  -0.19%  '

# Row 14
Set-TextValue 'D14' '36.26'
Set-TextValue 'E14' '  -2.29%  '

# Row 15
Set-TextValue 'D15' '4.381.84'
Set-TextValue 'E15' '  -2.08%  '

# Row 16
Set-TextValue 'D16' '3.742.92'
Set-TextValue 'E16' '  -2.38%  '

# Row 17
Set-TextValue 'D17' '18.55'
Set-TextValue 'E17' '  +0.15%  '

# Row 18
Set-TextValue 'D18' '67.375.59'
Set-TextValue 'E18' '  -1.40%  '

# Row 19
Set-TextValue 'D19' '7.15'
Set-TextValue 'E19' '  -3.17%  '

# Row 20
Set-TextValue 'E20' '  +0.71%  '

# Row 21
Set-TextValue 'D21' '10.48'
Set-TextValue 'E21' '  -5.46%  '

# Row 22
Set-TextValue 'D22' '466.30'
Set-TextValue 'E22' '  -0.63%  '

# Row 23
Set-TextValue 'D23' '0.716'
Set-TextValue 'E23' '  -2.61%  '

# Row 24
Set-TextValue 'D24' '83.44'
Set-TextValue 'E24' '  +0.34%  '

# Row 25
Set-TextValue 'E25' '  -8.74%  '

# Row 26
Set-TextValue 'D26' '2.19'
Set-TextValue 'E26' '  -1.74%  '

# Row 27
Set-TextValue 'D27' '12.11'
Set-TextValue 'E27' '  -0.45%  '

# Row 28
Set-TextValue 'D28' '10.24'
Set-TextValue 'E28' '  +1.77%  '

# Row 29
Set-TextValue 'E29' '  +0.04%  '

# Row 30
Set-TextValue 'D30' '2.89'
Set-TextValue 'E30' '  -2.64%  '

# Row 31
Set-TextValue 'D31' '3.900.82'
Set-TextValue 'E31' '  -1.98%  '

# Row 32
Set-TextValue 'D32' '7.61'
Set-TextValue 'E32' '  -1.22%  '

# Row 33
Set-TextValue 'D33' '30.31'
Set-TextValue 'E33' '  -4.12%  '

# Row 34
Set-TextValue 'E34' '  -4.01%  '

# Row 35
Set-TextValue 'D35' '9.08'
Set-TextValue 'E35' '  -3.32%  '

# Row 36
Set-TextValue 'D36' '3.711.71'
Set-TextValue 'E36' '  -2.25%  '

# Row 37
Set-TextValue 'D37' '3.78'
Set-TextValue 'E37' '  +1.95%  '

# Row 38
Set-TextValue 'E38' '  -1.43%  '

# Row 39
Set-TextValue 'D39' '0.138'
Set-TextValue 'E39' '  -1.73%  '

# Row 40
Set-TextValue 'D40' '0.997'
Set-TextValue 'E40' '  -2.05%  '

# Row 41
Set-TextValue 'D41' '5.78'
Set-TextValue 'E41' '  -2.65%  '

# Row 42
Set-TextValue 'E42' '  +0.07%  '

# Row 43
Set-TextValue 'D43' '0.310'
Set-TextValue 'E43' '  -1.47%  '

# Row 45
Set-TextValue 'D45' '8.66'
Set-TextValue 'E45' '  -0.98%  '

# Row 46
Set-TextValue 'D46' '1.93'
Set-TextValue 'E46' '  -2.44%  '

# Row 47
Set-TextValue 'E47' '  -2.76%  '

# Row 48
Set-TextValue 'D48' '395.48'
Set-TextValue 'E48' '  -5.04%  '

# Row 49
Set-TextValue 'D49' '0.000269'
Set-TextValue 'E49' '  -6.88%  '

# Row 50
Set-TextValue 'B50' 'Monero'
Set-TextValue 'C50' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D50' '138.77'
Set-TextValue 'E50' '  -2.12%  '

# Row 51
Set-TextValue 'B51' 'VeChain'
Set-TextValue 'C51' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D51' '0.0352'
Set-TextValue 'E51' '  -2.30%  '
